# Fix the erroneous option text in F18 ("both b & c" -> correct formula text),
# then remove the now fully-duplicate row 68 (same question as row 18) that had
# the correct text all along. Deleting the row shifts rows 69-71 up to 68-70.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1. Fix the broken answer-option text on row 18 ("The power can be expressed as").
$ws.Range("F18").Value = "both <P = VI> and <P = I^2R>"

# 2. Row 68 duplicated row 18's question but (unlike row 18) already had the
#    correct option text. Now that row 18 is fixed, row 68 is fully redundant,
#    so delete it entirely (rows below shift up).
$ws.Rows("68:68").Delete()
